# Apply cryptocurrency price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while forcing text storage so that
# numeric-looking strings (e.g. "253.72") are not silently converted
# into real numbers by Excel's auto-detection.
function Set-TextCell($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value2 = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" "42.881.85"
$ws.Range("E2").Value2 = "  +3.54%  "

Set-TextCell "D3" "2.255.16"
$ws.Range("E3").Value2 = "  +2.91%  "

$ws.Range("E4").Value2 = "  +0.04%  "

Set-TextCell "D5" "253.72"
$ws.Range("E5").Value2 = "  -0.61%  "

Set-TextCell "D6" "0.625"
$ws.Range("E6").Value2 = "  -0.59%  "

Set-TextCell "D7" "72.06"
$ws.Range("E7").Value2 = "  +4.86%  "

$ws.Range("E8").Value2 = "  -0.11%  "

Set-TextCell "D9" "0.649"
$ws.Range("E9").Value2 = "  +12.09%  "

Set-TextCell "D10" "41.11"
$ws.Range("E10").Value2 = "  +8.14%  "

Set-TextCell "D11" "59.56"
$ws.Range("E11").Value2 = "  +0.86%  "

Set-TextCell "D12" "0.0968"
$ws.Range("E12").Value2 = "  +3.10%  "

$ws.Range("E13").Value2 = "  +3.37%  "

$ws.Range("E14").Value2 = "  +0.71%  "

Set-TextCell "D15" "2.595.34"
$ws.Range("E15").Value2 = "  +3.21%  "

Set-TextCell "D16" "0.885"
$ws.Range("E16").Value2 = "  +0.87%  "

Set-TextCell "D17" "14.75"
$ws.Range("E17").Value2 = "  +1.36%  "

Set-TextCell "D18" "2.258.82"
$ws.Range("E18").Value2 = "  +1.49%  "

Set-TextCell "D19" "42.815.21"
$ws.Range("E19").Value2 = "  +3.66%  "

$ws.Range("E20").Value2 = "  +1.86%  "

$ws.Range("E21").Value2 = "  +0.90%  "

Set-TextCell "D22" "73.15"
$ws.Range("E22").Value2 = "  +1.45%  "

Set-TextCell "D23" "236.18"
$ws.Range("E23").Value2 = "  +1.28%  "

$ws.Range("E24").Value2 = "  +4.19%  "

Set-TextCell "D25" "3.97"
$ws.Range("E25").Value2 = "  +0.10%  "

$ws.Range("E26").Value2 = "  -1.08%  "

$ws.Range("E27").Value2 = "  +0.10%  "

Set-TextCell "D28" "2.45"
$ws.Range("E28").Value2 = "  -2.94%  "

Set-TextCell "D29" "3.69"
$ws.Range("E29").Value2 = "  -0.69%  "

$ws.Range("E30").Value2 = "  +2.12%  "

Set-TextCell "D31" "167.88"
$ws.Range("E31").Value2 = "  -0.55%  "

Set-TextCell "D32" "21.04"
$ws.Range("E32").Value2 = "  +1.62%  "

Set-TextCell "D33" "0.128"
$ws.Range("E33").Value2 = "  +8.80%  "

Set-TextCell "D34" "6.16"
$ws.Range("E34").Value2 = "  +11.99%  "

Set-TextCell "D35" "0.0785"
$ws.Range("E35").Value2 = "  +3.63%  "

$ws.Range("E36").Value2 = "  +1.17%  "

Set-TextCell "D37" "28.94"
$ws.Range("E37").Value2 = "  +7.90%  "

Set-TextCell "D38" "4.73"
$ws.Range("E38").Value2 = "  +1.97%  "

Set-TextCell "D39" "4.11"
$ws.Range("E39").Value2 = "  -1.21%  "

Set-TextCell "D40" "0.0321"
$ws.Range("E40").Value2 = "  +7.37%  "

Set-TextCell "D41" "2.29"
$ws.Range("E41").Value2 = "  +3.85%  "

Set-TextCell "D42" "5.94"
$ws.Range("E42").Value2 = "  +4.23%  "

Set-TextCell "D43" "12.40"
$ws.Range("E43").Value2 = "  -0.72%  "

Set-TextCell "D44" "64.26"
$ws.Range("E44").Value2 = "  -0.44%  "

Set-TextCell "D45" "4.98"
$ws.Range("E45").Value2 = "  -2.83%  "

Set-TextCell "D46" "9.06"
$ws.Range("E46").Value2 = "  +4.82%  "

$ws.Range("E47").Value2 = "  -0.38%  "

$ws.Range("E48").Value2 = "  +0.60%  "

$ws.Range("E49").Value2 = "  +0.88%  "

$ws.Range("E50").Value2 = "  -0.13%  "

$ws.Range("B51").Value2 = "TrustWalletToken"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D51" "1.19"
$ws.Range("E51").Value2 = "  +1.10%  "
